$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the three existing rows whose open/high/low/close values were revised ---
$ws.Cells.Item(253, 3).Value = 5817095330000
$ws.Cells.Item(253, 4).Value = 5817095330000
$ws.Cells.Item(253, 5).Value = 5817095330000
$ws.Cells.Item(253, 6).Value = 5817095330000

$ws.Cells.Item(254, 3).Value = 5950864520000
$ws.Cells.Item(254, 4).Value = 5950864520000
$ws.Cells.Item(254, 5).Value = 5950864520000
$ws.Cells.Item(254, 6).Value = 5950864520000

$ws.Cells.Item(255, 3).Value = 6045092150000
$ws.Cells.Item(255, 4).Value = 6045092150000
$ws.Cells.Item(255, 5).Value = 6045092150000
$ws.Cells.Item(255, 6).Value = 6045092150000

# --- Append three new monthly data rows (256-258) ---
$dateFormat = $ws.Cells.Item(255, 1).NumberFormat

$newRows = @(
    @{ Row = 256; Date = 44986.45833333334; Value = 6077620130000 },
    @{ Row = 257; Date = 45017.45833333334; Value = 6141246740000 },
    @{ Row = 258; Date = 45047.41666666666; Value = 6224272840000 }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value = $item.Date
    $dateCell.NumberFormat = $dateFormat

    $ws.Cells.Item($r, 2).Value = "ECONOMICS:CZM2"
    $ws.Cells.Item($r, 3).Value = $item.Value
    $ws.Cells.Item($r, 4).Value = $item.Value
    $ws.Cells.Item($r, 5).Value = $item.Value
    $ws.Cells.Item($r, 6).Value = $item.Value
    $ws.Cells.Item($r, 7).Value = 0
}
